# Auto-generated edit script: update crypto price/volume columns (D, E)
# per commit "Updated cryptos list on Fri Jan 26 15:16:25 UTC 2024 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "41.472.00"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +3.80%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.256.04"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +2.21%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "302.73"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.04%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "91.83"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +4.75%  "
$ws.Range("E7").Value = "  +2.27%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.00"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.03%  "
$ws.Range("E9").Value = "  +3.51%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "53.98"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +8.39%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "32.18"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +7.46%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0795"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +2.40%  "
$ws.Range("E13").Value = "  +3.03%  "
$ws.Range("E14").Value = "  +2.98%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.602.31"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.99%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.16"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +3.54%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.188.90"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.73%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.751"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +3.67%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "41.359.92"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +3.73%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.21"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +8.84%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0₃0905"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.56%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.91"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.67%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "66.88"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.64%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "240.78"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.65%  "
$ws.Range("E25").Value = "  +4.87%  "
$ws.Range("E26").Value = "  -0.20%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.87"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +3.39%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "23.77"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +5.78%  "
$ws.Range("E29").Value = "  +6.08%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "9.65"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +5.58%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "158.02"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.67%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "33.74"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +8.10%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.00"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.09%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.19"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +6.37%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0737"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +4.03%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.05"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +7.85%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.37"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.50%  "
$ws.Range("E38").Value = "  +9.50%  "
$ws.Range("E39").Value = "  +2.84%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.103"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +6.09%  "
$ws.Range("E41").Value = "  +6.09%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.98"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +7.07%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "20.27"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +17.06%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.063.20"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.55%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0277"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +3.86%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "10.16"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +4.85%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.98"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +12.50%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.08"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.69%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.474.48"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.15%  "
$ws.Range("E50").Value = "  +2.98%  "
$ws.Range("E51").Value = "  +3.67%  "
